# Add a new worksheet "XMOS Multichip Planning" after "XMOS Multichip" to help
# plan XCore division of labor, and move the active/selected tab to it.

$wb = $excel.ActiveWorkbook

$xmosMultichip = $wb.Worksheets.Item("XMOS Multichip")

$planning = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $xmosMultichip)
$planning.Name = "XMOS Multichip Planning"

# Column widths (closest achievable given COM's pixel-grid quantization of
# ColumnWidth; targets are ~13.14 and ~21.43 "characters").
$planning.Columns.Item(1).ColumnWidth = 12.25
$planning.Columns.Item(3).ColumnWidth = 20.59

# Header row
$planning.Range("A1").Value = "Block"
$planning.Range("B1").Value = "1-bit ports"
$planning.Range("D1").Value = "Xcore"

# Data rows
$planning.Range("A2").Value = "Pmod 0"
$planning.Range("B2").Value = 8
$planning.Range("D2").Value = "U3"

$planning.Range("A3").Value = "Pmod 1"
$planning.Range("B3").Value = 8
$planning.Range("D3").Value = "U4"

$planning.Range("A4").Value = "Pmod 2"
$planning.Range("B4").Value = 8
$planning.Range("D4").Value = "U4"

$planning.Range("A5").Value = "Gadgeteer"
$planning.Range("B5").Value = 7
$planning.Range("D5").Value = "U3"

$planning.Range("A6").Value = "SPI"
$planning.Range("B6").Value = 5
$planning.Range("D6").Value = "U1"

$planning.Range("A7").Value = "I2C"
$planning.Range("B7").Value = 3
$planning.Range("C7").Value = "and a P4C"
$planning.Range("D7").Value = "U2"

$planning.Range("A8").Value = "PWM"
$planning.Range("B8").Value = 9
$planning.Range("D8").Value = "U1"

$planning.Range("A9").Value = "Enc"
$planning.Range("B9").Value = 6
$planning.Range("D9").Value = "U2"

$planning.Range("A10").Value = "STM32 UART"
$planning.Range("B10").Value = 3
$planning.Range("D10").Value = "U2"

$planning.Range("A11").Value = "Xbee UART"
$planning.Range("B11").Value = 2
$planning.Range("C11").Value = "2 free for RTS/CTS"
$planning.Range("D11").Value = "U2"

# Match the authored selection/view on the new (now active) sheet.
$planning.Range("C17").Select()

# The previously active sheet scrolls so row 52 is at the top; its tab is no
# longer the selected one (handled automatically by activating $planning above).
$xmosMultichip.Activate()
$excel.ActiveWindow.ScrollRow = 52
$excel.ActiveWindow.ScrollColumn = 1
$planning.Activate()
